# The "End Date" value in B6 (previously the serial date 43465 / 12/31/2018,
# shown with a short m/d/yyyy number format) is replaced with a literal text
# date string, including the source's own typo ("Decmeber"). The cell's
# number format is updated to a long-date custom format as part of the same
# edit (it stays registered on the cell even though the value is now text).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").NumberFormat = "[$-409]mmmm\ d\,\ yyyy;@"
$ws.Range("B6").Value = "Decmeber 31, 2018"

# The author's cursor ended up on B11 when the file was last saved.
$ws.Range("B11").Select()
